$wb = $excel.ActiveWorkbook

# 展览 (was sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 27086
$ws.Range("F4").Value = 671
$ws.Range("F5").Value = 202
$ws.Range("F7").Value = 235
$ws.Range("F9").Value = 500
$ws.Range("F10").Value = 201
$ws.Range("F13").Value = 107
$ws.Range("F14").Value = 525
$ws.Range("F16").Value = 1665
$ws.Range("F17").Value = 277
$ws.Range("F18").Value = 1130
$ws.Range("F19").Value = 204
$ws.Range("F20").Value = 469
$ws.Range("F21").Value = 15
$ws.Range("F23").Value = 124

# 演出 (was sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 4534
$ws.Range("F10").Value = 460

# 本地生活 (was sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5209
$ws.Range("F3").Value = 284

# 全部类型 (was sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5209
$ws.Range("F4").Value = 284
$ws.Range("F5").Value = 27086
$ws.Range("F6").Value = 4534
$ws.Range("F7").Value = 671
$ws.Range("F10").Value = 202
$ws.Range("F16").Value = 460
$ws.Range("F20").Value = 235
$ws.Range("F22").Value = 500
$ws.Range("F23").Value = 201
$ws.Range("F27").Value = 107
$ws.Range("F30").Value = 525
$ws.Range("F33").Value = 1665
$ws.Range("F34").Value = 277
$ws.Range("F35").Value = 1130
$ws.Range("F37").Value = 204
$ws.Range("F38").Value = 469
$ws.Range("F39").Value = 15
$ws.Range("F42").Value = 124
